$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (regenerated to use K instead of Strike#).
# Update rows 2-10 with the newly calculated values.
$ws.Range("G2").Value = 7
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 4
$ws.Range("G10").Value = 2
